$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 257.5
$ws.Range("I12").Value = 257.5
$ws.Range("K12").Value = 257.5
$ws.Range("M12").Value = -87.5
$ws.Range("H132").Value = 2800.6875
$ws.Range("I132").Value = 3178.8076
$ws.Range("J132").Value = 1162.1666
$ws.Range("K132").Value = 9536.4228
$ws.Range("L132").Value = 3486.4998
$ws.Range("M132").Value = -7006.4228
$ws.Range("N132").Value = -8546.4998
$ws.Range("H137").Value = 395929.75
$ws.Range("I137").Value = 2058.75
$ws.Range("K137").Value = 6176.25
$ws.Range("M137").Value = -3626.25
$ws.Range("H141").Value = 3629.4443
$ws.Range("I141").Value = 622.06665
$ws.Range("K141").Value = 1866.19995
$ws.Range("M141").Value = 3313.80005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4249.04
$ws.Range("I32").Value = 3488.8923
$ws.Range("J32").Value = 9190
$ws.Range("K32").Value = 3488.8923
$ws.Range("L32").Value = 9190
$ws.Range("M32").Value = -3201.8923
$ws.Range("N32").Value = -9764
$ws.Range("H61").Value = 4717.9688
$ws.Range("I61").Value = 2560.8
$ws.Range("K61").Value = 2560.8
$ws.Range("M61").Value = -2348.8
$ws.Range("H63").Value = 3005
$ws.Range("I63").Value = 3005
$ws.Range("K63").Value = 3005
$ws.Range("M63").Value = -2319
$ws.Range("H66").Value = 3005
$ws.Range("I66").Value = 3005
$ws.Range("K66").Value = 15025
$ws.Range("M66").Value = -11593
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H74").Value = 835415.8
$ws.Range("I74").Value = 1668581.6
$ws.Range("K74").Value = 1668581.6
$ws.Range("M74").Value = -1667707.6
$ws.Range("H77").Value = 835415.8
$ws.Range("I77").Value = 1668581.6
$ws.Range("K77").Value = 8342908
$ws.Range("M77").Value = -8338540
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H113").Value = 58331
$ws.Range("J113").Value = 58331
$ws.Range("L113").Value = 58331
$ws.Range("N113").Value = -67009
$ws.Range("H136").Value = 4717.9688
$ws.Range("I136").Value = 2560.8
$ws.Range("K136").Value = 7682.400000000001
$ws.Range("M136").Value = -5132.400000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1102.2222
$ws.Range("I7").Value = 1225
$ws.Range("J7").Value = 1004
$ws.Range("K7").Value = 1225
$ws.Range("L7").Value = 1004
$ws.Range("M7").Value = -1112
$ws.Range("N7").Value = -1230
$ws.Range("H20").Value = 1778.8572
$ws.Range("I20").Value = 1936.1818
$ws.Range("K20").Value = 1936.1818
$ws.Range("M20").Value = -1689.1818
$ws.Range("H93").Value = 15000
$ws.Range("I93").Value = 15000
$ws.Range("K93").Value = 15000
$ws.Range("M93").Value = -13128
$ws.Range("H94").Value = 26974.766
$ws.Range("I94").Value = 342.29166
$ws.Range("J94").Value = 90892.7
$ws.Range("K94").Value = 342.29166
$ws.Range("L94").Value = 90892.7
$ws.Range("M94").Value = 108.70834
$ws.Range("N94").Value = -91794.7
$ws.Range("H95").Value = 17285.5
$ws.Range("J95").Value = 17285.5
$ws.Range("L95").Value = 17285.5
$ws.Range("N95").Value = -22777.5
$ws.Range("H96").Value = 8428
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H97").Value = 8914.25
$ws.Range("I97").Value = 818.4
$ws.Range("J97").Value = 14697
$ws.Range("K97").Value = 818.4
$ws.Range("L97").Value = 14697
$ws.Range("M97").Value = 172.6
$ws.Range("N97").Value = -16679
$ws.Range("H98").Value = 60542
$ws.Range("J98").Value = 60542
$ws.Range("L98").Value = 60542
$ws.Range("N98").Value = -66532
$ws.Range("H100").Value = 19403.857
$ws.Range("J100").Value = 19403.857
$ws.Range("L100").Value = 19403.857
$ws.Range("N100").Value = -21567.857
$ws.Range("H102").Value = 9709.5
$ws.Range("I102").Value = 3651.4
$ws.Range("K102").Value = 3651.4
$ws.Range("M102").Value = -406.4000000000001
$ws.Range("H103").Value = 16111.846
$ws.Range("J103").Value = 16111.846
$ws.Range("L103").Value = 16111.846
$ws.Range("N103").Value = -18455.846

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2063.03
$ws.Range("I31").Value = 1647.75
$ws.Range("J31").Value = 2512.9167
$ws.Range("K31").Value = 1647.75
$ws.Range("L31").Value = 2512.9167
$ws.Range("M31").Value = -1352.75
$ws.Range("N31").Value = -3102.9167
$ws.Range("H34").Value = 2063.03
$ws.Range("I34").Value = 1647.75
$ws.Range("J34").Value = 2512.9167
$ws.Range("K34").Value = 1647.75
$ws.Range("L34").Value = 2512.9167
$ws.Range("M34").Value = -1445.75
$ws.Range("N34").Value = -2916.9167
$ws.Range("H58").Value = 2437.9429
$ws.Range("I58").Value = 1974.7
$ws.Range("J58").Value = 5217.4
$ws.Range("K58").Value = 1974.7
$ws.Range("L58").Value = 5217.4
$ws.Range("M58").Value = -1771.7
$ws.Range("N58").Value = -5623.4
$ws.Range("H134").Value = 1848.4595
$ws.Range("I134").Value = 1665.5143
$ws.Range("J134").Value = 5050
$ws.Range("K134").Value = 4996.5429
$ws.Range("L134").Value = 15150
$ws.Range("M134").Value = -2461.5429
$ws.Range("N134").Value = -20220
$ws.Range("H136").Value = 2437.9429
$ws.Range("I136").Value = 1974.7
$ws.Range("J136").Value = 5217.4
$ws.Range("K136").Value = 5924.1
$ws.Range("L136").Value = 15652.2
$ws.Range("M136").Value = -3374.1
$ws.Range("N136").Value = -20752.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 685.7826
$ws.Range("I5").Value = 551.0952
$ws.Range("K5").Value = 1653.2856
$ws.Range("M5").Value = -1541.2856
$ws.Range("H107").Value = 1708.1875
$ws.Range("J107").Value = 2295.7273
$ws.Range("L107").Value = 6887.1819
$ws.Range("N107").Value = -10727.1819
$ws.Range("H122").Value = 1196.7
$ws.Range("J122").Value = 1458.5
$ws.Range("L122").Value = 13126.5
$ws.Range("N122").Value = -18026.5
$ws.Range("H131").Value = 3076
$ws.Range("J131").Value = 5986.2
$ws.Range("L131").Value = 17958.6
$ws.Range("N131").Value = -28038.6
$ws.Range("H135").Value = 685.7826
$ws.Range("I135").Value = 551.0952
$ws.Range("K135").Value = 4959.8568
$ws.Range("M135").Value = -2424.8568
$ws.Range("H137").Value = 3785.8765
$ws.Range("I137").Value = 1297.125
$ws.Range("J137").Value = 4031.679
$ws.Range("K137").Value = 3891.375
$ws.Range("L137").Value = 12095.037
$ws.Range("M137").Value = 1208.625
$ws.Range("N137").Value = -22295.037

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2893.45
$ws.Range("I132").Value = 4098.6665
$ws.Range("J132").Value = 2170.32
$ws.Range("K132").Value = 12295.9995
$ws.Range("L132").Value = 6510.960000000001
$ws.Range("M132").Value = -9765.999500000002
$ws.Range("N132").Value = -11570.96

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1629.7
$ws.Range("I22").Value = 924.75
$ws.Range("K22").Value = 924.75
$ws.Range("M22").Value = -629.75
$ws.Range("H27").Value = 1629.7
$ws.Range("I27").Value = 924.75
$ws.Range("K27").Value = 924.75
$ws.Range("M27").Value = -817.75
$ws.Range("H55").Value = 1290.25
$ws.Range("J55").Value = 2806.5833
$ws.Range("L55").Value = 2806.5833
$ws.Range("N55").Value = -3152.5833
$ws.Range("H114").Value = 60000
$ws.Range("J114").Value = 60000
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678
$ws.Range("H122").Value = 3189.2886
$ws.Range("I122").Value = 3221.8958
$ws.Range("J122").Value = 2798
$ws.Range("K122").Value = 9665.687399999999
$ws.Range("L122").Value = 8394
$ws.Range("M122").Value = -7215.687399999999
$ws.Range("N122").Value = -13294

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1079.091
$ws.Range("I81").Value = 1079.091
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2158.182
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1097.182
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 1079.091
$ws.Range("I84").Value = 1079.091
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10790.91
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -5486.91
$ws.Range("N84").Value = $null
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = $null
$ws.Range("H112").Value = 18907.334
$ws.Range("J112").Value = 18907.334
$ws.Range("L112").Value = 18907.334
$ws.Range("N112").Value = -21861.334
$ws.Range("H121").Value = 57209.5
$ws.Range("J121").Value = 57209.5
$ws.Range("L121").Value = 57209.5
$ws.Range("N121").Value = -60703.5
$ws.Range("H124").Value = 23970
$ws.Range("J124").Value = 23970
$ws.Range("L124").Value = 23970
$ws.Range("N124").Value = -33790
